# "update pp19 barrel and weight"
# Adds a new PP19 barrel mod (row 19) to the m4-barrels sheet, wires the
# N (balance) and S (irl price check) shared-formula columns down through
# the new row, fills in row 18 (which was previously a blank spacer row but
# now participates in the N formula range), and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: give it the balance formula (it's now inside the N4:N19 run) ---
$ws.Range("N18").Formula = "=C18-D18*20-E18*0.8-F18*0.6-H18*5+I18*10+J18/300"

# --- Row 19: new barrel entry ---
$ws.Range("A19").Value = "pp19_barrel"
$ws.Range("B19").Value = "PP19 Standard"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0.3
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("H19").Value = 0.1
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 100
$ws.Range("M19").Value = 0

$ws.Range("N19").Formula = "=C19-D19*20-E19*0.8-F19*0.6-H19*5+I19*10+J19/300"

$ws.Range("P19").Value = 0.06
$ws.Range("Q19").Value = 9

$ws.Range("S19").Formula = "=ROUND(Q19*0.02+P19+R19, 2)"

# --- move the selection cursor like the author's last click ---
$ws.Range("E24").Select() | Out-Null
